$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 should become the "pv1"/"pv" entry (it currently holds demand2/demand)
$ws.Range("A3").Value = "pv1"
$ws.Range("B3").Value = "pv"

# Delete rows 4 through 19, keeping only header (row1), demand1 (row2) and pv1 (row3)
$ws.Range("A4:B19").EntireRow.Delete()
